# === Sheet1 updates ===
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Convert CONCATENATE formulas in E2:E21 to HYPERLINK(CONCATENATE(...)) and apply Hyperlink style
$ws1.Range("E2").Formula = '=HYPERLINK(CONCATENATE("https://www.kaggle.com/",B2,"/",C2))'
$ws1.Range("E3:E21").Formula = '=HYPERLINK(CONCATENATE("https://www.kaggle.com/",B3,"/",C3))'
$ws1.Range("E2:E21").Style = "Hyperlink"

# Widen column D on Sheet1
$ws1.Columns.Item(4).ColumnWidth = 18.25

# New row 23 noting nbconvert limitation
$ws1.Range("B23").Value = "Nbconvert encoding support limited"

# === Sheet2: new worksheet with notebook-forking analysis notes ===
$ws2 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws2.Name = "Sheet2"

# Header row (no special style)
$ws2.Range("B1").Value = "Forked code"
$ws2.Range("C1").Value = "Changes"

# Data rows 2-12: forking category (col B) and observed changes (col C/D)
$ws2.Range("B2").Value = "Documentation"
$ws2.Range("C2").Value = "Add documentation"

$ws2.Range("B3").Value = "Data processing and splitting"
$ws2.Range("C3").Value = "Append model training code to data processing code (may benefit from forking tool due to the need for trial and error during training)."

$ws2.Range("B4").Value = "Simple ML model with data processing"
$ws2.Range("C4").Value = "New features (e.g., GPU support), model changes, different model types (may benefit from forking), different loss functions (may benefit from forking)"

$ws2.Range("B5").Value = "Simple ML model with data processing"
$ws2.Range("C5").Value = "Parameter changes, different model type, changes to data processing (due to the different model)"

$ws2.Range("B6").Value = "Data processing and splitting"
$ws2.Range("C6").Value = "Different weights, added a model to the starter notebook which only contained data processing, added data processing for training data and validation data"

$ws2.Range("B7").Value = "Tool usage template"
$ws2.Range("C7").Value = "Adapted original notebook to add enhanced data processing. Original notebook served as a baseline for basic usage of a tool. Original data did not change"

$ws2.Range("B8").Value = "Model evaluation and comparison (from model files); data processing"
$ws2.Range("C8").Value = "Added plotting of results, changed most of the data processing, Adapted original notebook for inference problem using a model file instead of training and model accuracy evaluation of a set of models. Author has made other notebooks which fork from same origin as this one."
$ws2.Range("D8").Value = "Q: how to fix this? A lot repeated code, but also a lot that is different"

$ws2.Range("B9").Value = "Data processing and exploration"
$ws2.Range("C9").Value = "Starts diverging after the 5th cell (forking could be useful); added train test split; data sampling is slightly different; replace data visualization and data output portions with model training - added model training, model selection, model evaluation; "

$ws2.Range("B10").Value = "Simple ML model with data processing"
$ws2.Range("C10").Value = "Removed some feature engineering code, different plots in some cases, small changes to the way relations between variables are displayed. Model training section did not change"

$ws2.Range("B11").Value = "Simple ML model with data processing"
$ws2.Range("C11").Value = "small variable and parameter changes during data processing; large changes in the model. This notebook was used for an alternative exploration; can incorporate forking"

$ws2.Range("B12").Value = "Demonstration of a library"
$ws2.Range("C12").Value = "Removed and added import statements; added additional functions; code is different everywhere except for the majority of the import statements"

# Row index column A (0-19) for rows 2-21, reuse the bordered/bold style from Sheet1 col A
$ws1.Range("A2:A21").Copy()
$ws2.Range("A2:A21").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false
For ($i = 0; $i -le 19; $i++) {
  $ws2.Cells.Item($i + 2, 1).Value = $i
}

# Apply wrap-text style to notes columns B:D for rows 2-21
$ws2.Range("B2:D21").WrapText = $true

# Column widths
$ws2.Columns.Item(2).ColumnWidth = 93.59
$ws2.Columns.Item(3).ColumnWidth = 71.42
$ws2.Columns.Item(4).ColumnWidth = 50.25

# Row heights for wrapped rows (matches rendered text wrapping)
$ws2.Rows.Item(3).RowHeight = 30
$ws2.Rows.Item(4).RowHeight = 30
$ws2.Rows.Item(5).RowHeight = 30
$ws2.Rows.Item(6).RowHeight = 45
$ws2.Rows.Item(7).RowHeight = 45
$ws2.Rows.Item(8).RowHeight = 60
$ws2.Rows.Item(9).RowHeight = 60
$ws2.Rows.Item(10).RowHeight = 45
$ws2.Rows.Item(11).RowHeight = 45
$ws2.Rows.Item(12).RowHeight = 30

# Footer notes: general forking ideas
$ws2.Range("B25").Value = "Idea: forking between notebooks and smart merge or merge based on user defined cells, benefits original code"
$ws2.Range("B26").Value = "Idea 2: fork from multiple sources (pick and choose)"

# Final selections and active sheet/window state
$ws1.Range("C9").Select()
$ws2.Range("B13").Select()
$ws2.Activate()
